$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 33, shifting existing rows 33+ down by one.
$ws.Rows.Item(33).Insert()

# Populate the newly inserted row 33 with the new weekly record.
$ws.Cells.Item(33, 1).Value = 9
$ws.Cells.Item(33, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(33, 3).Value = "Metropolitana"
$ws.Cells.Item(33, 4).Value = 44714
$ws.Cells.Item(33, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(33, 5).Value = 13
$ws.Cells.Item(33, 6).Value = 100112035
$ws.Cells.Item(33, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(33, 8).Value = "Sin especificar"
$ws.Cells.Item(33, 9).Value = "Primera"
$ws.Cells.Item(33, 10).Value = 52
$ws.Cells.Item(33, 11).Value = 18000
$ws.Cells.Item(33, 12).Value = 20000
$ws.Cells.Item(33, 13).Value = 19000
$ws.Cells.Item(33, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(33, 15).Value = "Hijuelas"
$ws.Cells.Item(33, 16).Value = 1267
$ws.Cells.Item(33, 17).Value = 15
$ws.Cells.Item(33, 18).Value = "Hortaliza"
